$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains price figures that look numeric (e.g. "1.001", "0.08490")
# but must be preserved as literal text (matching the source inline strings),
# including significant trailing zeros and multi-dot "thousands.decimal" values.
# Force the column to Text format first so assignment does not coerce to Double.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.230.47"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "1.647.32"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "308.89"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.3909"
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("D8").Value = "0.3879"
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").Value = "1.001"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "1.369"
$ws.Range("E10").Value = "  -6.83%  "
$ws.Range("D11").Value = "48.66"
$ws.Range("E11").Value = "  -7.92%  "
$ws.Range("D12").Value = "0.08490"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").Value = "24.31"
$ws.Range("E13").Value = "  -6.36%  "
$ws.Range("D14").Value = "7.197"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").Value = "0.00001292"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").Value = "7.550"
$ws.Range("E16").Value = "  -5.21%  "
$ws.Range("D17").Value = "1.650.61"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").Value = "95.07"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "0.06951"
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("D20").Value = "21.20"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("D21").Value = "6.985"
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "13.85"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "24.224.05"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "2.347"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "2.737"
$ws.Range("E26").Value = "  -6.71%  "
$ws.Range("D27").Value = "22.62"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("D28").Value = "8.959"
$ws.Range("E28").Value = "  +7.87%  "
$ws.Range("D29").Value = "158.15"
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").Value = "142.56"
$ws.Range("E30").Value = "  -4.89%  "
$ws.Range("D31").Value = "5.381"
$ws.Range("E31").Value = "  -12.44%  "
$ws.Range("D32").Value = "2.449"
$ws.Range("E32").Value = "  -7.58%  "
$ws.Range("D33").Value = "1.829.37"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("D34").Value = "7.234"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").Value = "0.08116"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("D36").Value = "0.9923"
$ws.Range("E36").Value = "  -4.66%  "
$ws.Range("D37").Value = "0.02960"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").Value = "0.2724"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("D39").Value = "0.09307"
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("D40").Value = "1.483"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "10.11"
$ws.Range("E41").Value = "  -7.12%  "
$ws.Range("D42").Value = "0.7680"
$ws.Range("E42").Value = "  -6.82%  "
$ws.Range("D43").Value = "13.17"
$ws.Range("E43").Value = "  -5.71%  "
$ws.Range("D44").Value = "16.04"
$ws.Range("E44").Value = "  -6.40%  "
$ws.Range("D45").Value = "2.505"
$ws.Range("E45").Value = "  -6.52%  "
$ws.Range("D46").Value = "0.6925"
$ws.Range("E46").Value = "  -6.22%  "
$ws.Range("D47").Value = "4.104"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "0.08427"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "1.272"
$ws.Range("E50").Value = "  -9.40%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "134.54"
$ws.Range("E51").Value = "  -3.35%  "
